$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new column S that duplicates column J ("Sound Spd. (m/s)") -
# header plus the 11 data rows (rows 2-12).
$ws.Range("S1").Value = $ws.Range("J1").Value2
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells($r, 19).Value = $ws.Cells($r, 10).Value2
}

# Approximate the column widths Excel auto-derived for the newly
# touched columns R (18) and S (19).
$ws.Columns.Item(18).ColumnWidth = 14.6666666666667
$ws.Columns.Item(19).ColumnWidth = 15.5

# Reflect that column S got selected last (whole-column selection),
# matching the workbook's last recorded UI state.
$ws.Columns.Item(19).Select() | Out-Null
